# "Generate Report for Handback" — marks the zh-cn and de-de handback rows
# as complete: status flips to "Handed back: in sync with en-US", the
# Latest Target File / Latest Handback File / Latest Handback DateTime
# columns get filled in (with a hyperlink on the target-file cell), and a
# few columns are widened to fit the new, longer text.

$wb  = $excel.ActiveWorkbook
$ovw = $wb.Worksheets.Item("Overview")
$zh  = $wb.Worksheets.Item("zh-cn")
$de  = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5680811fc8899b4c0bd0069f4a18a689c58b2b3f/e2e/7860255e-61b4-4612-a2bb-63b1b2ec886e.md"
$urlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5680811fc8899b4c0bd0069f4a18a689c58b2b3f/e2e/81be7453-d1f2-4473-bc59-0b9cc878ac7e.md"
$nameA = "7860255e-61b4-4612-a2bb-63b1b2ec886e.md"
$nameB = "81be7453-d1f2-4473-bc59-0b9cc878ac7e.md"

# Color used by the workbook's existing custom "HyperLink" cell style
# (font rgb FF6495ED, underlined) so the newly-hyperlinked cells match the
# look of the already-hyperlinked "Source File Name" column.
$hyperlinkColor = 15570276

# ---------------------------------------------------------------------
# Overview sheet: both locale-status columns (E = zh-cn, F = de-de) show
# the same "Ready for handoff" text today; both move to the new status.
# ---------------------------------------------------------------------
$ovw.Range("E2").Value = $statusText
$ovw.Range("F2").Value = $statusText
$ovw.Range("E3").Value = $statusText
$ovw.Range("F3").Value = $statusText

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh.Range("C2").Value = $statusText
$zh.Range("C3").Value = $statusText

$zh.Hyperlinks.Add($zh.Range("I2"), $urlA, "", "", $nameA) | Out-Null
$zh.Range("I2").Font.Underline = $true
$zh.Range("I2").Font.Color = $hyperlinkColor

$zh.Hyperlinks.Add($zh.Range("I3"), $urlB, "", "", $nameB) | Out-Null
$zh.Range("I3").Font.Underline = $true
$zh.Range("I3").Font.Color = $hyperlinkColor

$zh.Range("J2").Value = "7860255e-61b4-4612-a2bb-63b1b2ec886e.a0f5e76e9cc141511fa94a51d1f44f0335f4ca33.zh-cn.xlf"
$zh.Range("J3").Value = "81be7453-d1f2-4473-bc59-0b9cc878ac7e.45b7a16c45c193b999653ed8f9ce75d7df741a0e.zh-cn.xlf"

$zh.Range("K2").Value = "2016-08-31 03:10:53"
$zh.Range("K3").Value = "2016-08-31 03:10:53"

# Columns C (Status), I (Latest Target File) and J (Latest Handback File)
# widen to fit the longer strings now in them.
$zh.Columns.Item(3).ColumnWidth  = 29.166666666666668
$zh.Columns.Item(9).ColumnWidth  = 39.166666666666664
$zh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de.Range("C2").Value = $statusText
$de.Range("C3").Value = $statusText

$de.Hyperlinks.Add($de.Range("I2"), $urlA, "", "", $nameA) | Out-Null
$de.Range("I2").Font.Underline = $true
$de.Range("I2").Font.Color = $hyperlinkColor

$de.Hyperlinks.Add($de.Range("I3"), $urlB, "", "", $nameB) | Out-Null
$de.Range("I3").Font.Underline = $true
$de.Range("I3").Font.Color = $hyperlinkColor

$de.Range("J2").Value = "7860255e-61b4-4612-a2bb-63b1b2ec886e.a0f5e76e9cc141511fa94a51d1f44f0335f4ca33.de-de.xlf"
$de.Range("J3").Value = "81be7453-d1f2-4473-bc59-0b9cc878ac7e.45b7a16c45c193b999653ed8f9ce75d7df741a0e.de-de.xlf"

$de.Range("K2").Value = "2016-08-31 03:11:01"
$de.Range("K3").Value = "2016-08-31 03:11:01"

$de.Columns.Item(3).ColumnWidth  = 29.166666666666668
$de.Columns.Item(9).ColumnWidth  = 39.166666666666664
$de.Columns.Item(10).ColumnWidth = 39.166666666666664
